$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Build a formula that evaluates to the literal string, write it,
    # then flatten via Copy + PasteSpecial(values) so the cell ends up
    # holding a plain text value (not a formula) without Excel
    # reinterpreting numeric-looking text as a Number and without any
    # style/number-format changes on the cell.
    $rng = $ws.Range($cellRef)
    $escaped = $text.Replace('"', '""')
    $rng.Formula = '="' + $escaped + '"'
    $rng.Copy()
    $rng.PasteSpecial(-4163)
}

# --- Coin/Link text columns (B/C) ---
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'

# --- Price (D) / Volume(1h) (E) columns ---
Set-TextValue 'D2' '25.910.54'
Set-TextValue 'E2' '  +0.29%  '
Set-TextValue 'D3' '1.635.62'
Set-TextValue 'E3' '  +0.09%  '
Set-TextValue 'E4' '  +0.23%  '
Set-TextValue 'D5' '214.60'
Set-TextValue 'E5' '  +0.01%  '
Set-TextValue 'E6' '  +1.11%  '
Set-TextValue 'E7' '  +0.24%  '
Set-TextValue 'E8' '  -0.79%  '
Set-TextValue 'E9' '  +0.48%  '
Set-TextValue 'D10' '19.61'
Set-TextValue 'E10' '  -0.15%  '
Set-TextValue 'D11' '0.0795'
Set-TextValue 'E11' '  +0.64%  '
Set-TextValue 'D12' '1.861.82'
Set-TextValue 'E12' '  +0.08%  '
Set-TextValue 'E13' '  -0.43%  '
Set-TextValue 'D14' '1.641.08'
Set-TextValue 'E14' '  +0.29%  '
Set-TextValue 'D15' '0.542'
Set-TextValue 'E15' '  -1.66%  '
Set-TextValue 'E16' '  -0.26%  '
Set-TextValue 'D17' '62.59'
Set-TextValue 'E17' '  -0.37%  '
Set-TextValue 'D18' '25.930.79'
Set-TextValue 'E18' '  +0.41%  '
Set-TextValue 'E19' '  +0.24%  '
Set-TextValue 'D20' '4.39'
Set-TextValue 'E20' '  -0.98%  '
Set-TextValue 'D21' '193.47'
Set-TextValue 'E21' '  +1.10%  '
Set-TextValue 'D22' '9.92'
Set-TextValue 'E22' '  -0.37%  '
Set-TextValue 'D23' '6.27'
Set-TextValue 'E23' '  -0.75%  '
Set-TextValue 'E24' '  +0.16%  '
Set-TextValue 'D25' '143.77'
Set-TextValue 'E25' '  +0.87%  '
Set-TextValue 'E26' '  +0.32%  '
Set-TextValue 'E27' '  +2.94%  '
Set-TextValue 'D28' '6.84'
Set-TextValue 'E28' '  +0.00%  '
Set-TextValue 'D29' '15.42'
Set-TextValue 'E29' '  -0.52%  '
Set-TextValue 'E30' '  +0.23%  '
Set-TextValue 'E31' '  +1.34%  '
Set-TextValue 'D32' '3.29'
Set-TextValue 'E32' '  -1.16%  '
Set-TextValue 'D33' '3.21'
Set-TextValue 'E33' '  -0.56%  '
Set-TextValue 'E34' '  -2.50%  '
Set-TextValue 'E35' '  +1.36%  '
Set-TextValue 'E36' '  -0.45%  '
Set-TextValue 'D37' '1.139.01'
Set-TextValue 'E37' '  -0.68%  '
Set-TextValue 'D38' '0.545'
Set-TextValue 'E38' '  -0.13%  '
Set-TextValue 'E39' '  -1.02%  '
Set-TextValue 'E40' '  +0.13%  '
Set-TextValue 'E41' '  +0.19%  '
Set-TextValue 'D42' '99.43'
Set-TextValue 'E42' '  -1.06%  '
Set-TextValue 'E43' '  -0.65%  '
Set-TextValue 'E44' '  -3.57%  '
Set-TextValue 'D45' '1.771.07'
Set-TextValue 'E45' '  +0.10%  '
Set-TextValue 'D46' '0.0₆0114'
Set-TextValue 'E46' '  +3.10%  '
Set-TextValue 'D47' '56.29'
Set-TextValue 'E47' '  +1.29%  '
Set-TextValue 'E48' '  +3.26%  '
Set-TextValue 'E49' '  -0.95%  '
Set-TextValue 'D50' '0.415'
Set-TextValue 'E50' '  -0.39%  '
Set-TextValue 'D51' '7.62'
Set-TextValue 'E51' '  +0.88%  '

$excel.CutCopyMode = $false
